$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("web 121")

# New column header "Q02" (quiz 02) in C1, mirroring B1 ("Q01")
$ws.Range("C1").Value = "Q02"

# Quiz 02 raw scores (out of 12) for each student, entered as formulas
# scaling the raw mark to a 0-10 scale, matching the existing B-column
# pattern used for quiz 01 (out of 10).
$ws.Range("C2").Formula = "=(6/12)*10"
$ws.Range("C3").Value = 0
$ws.Range("C4").Formula = "=(7/12)*10"
$ws.Range("C5").Formula = "=(8/12)*10"
$ws.Range("C6").Formula = "=(8/12)*10"
$ws.Range("C7").Formula = "=(10/12)*10"
$ws.Range("C8").Value = 0
$ws.Range("C9").Formula = "=(6/12)*10"
$ws.Range("C10").Formula = "=(7/12)*10"
$ws.Range("C11").Formula = "=(7/12)*10"
$ws.Range("C12").Formula = "=(5/12)*10"
$ws.Range("C13").Formula = "=(5/12)*10"
$ws.Range("C14").Formula = "=(5/12)*10"
$ws.Range("C15").Formula = "=(9/12)*10"

# Move the active selection to C16, matching the author's cursor position
$ws.Range("C16").Select()
